$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$np = $s.NotesPage
$tcs = $np.ThemeColorScheme
for ($i=1; $i -le 12; $i++) {
  Write-Output ("  $i=" + $tcs.Item($i).RGB)
}
$tcs.Item(1).RGB = 999999
Write-Output "set"
